$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 330, shifting existing data (rows 330-355) down to 332-357
$ws.Rows.Item(330).Resize(2).Insert()

# Common (constant) values for this block of rows
$mercadoId = 4
$mercado = "Feria Lagunitas de Puerto Montt"
$region = "Los Lagos"
$codreg = 10
$tipo = "Fruta"
$productoId = 100102
$producto = "Cítricos"
$categoriaId = 100102006
$categoria = "Pomelo"
$unidad = "`$/caja 14 kilos empedrada"
$origen = "Región de O'Higgins"
$kgUnidad = 14

# New row 330: Start Ruby / Primera
$r = 330
$ws.Cells.Item($r, 1).Value = $mercadoId
$ws.Cells.Item($r, 2).Value = $mercado
$ws.Cells.Item($r, 3).Value = $region
$ws.Cells.Item($r, 4).Value = 44783
$ws.Cells.Item($r, 4).NumberFormat = $ws.Cells.Item(332, 4).NumberFormat
$ws.Cells.Item($r, 5).Value = $codreg
$ws.Cells.Item($r, 6).Value = $tipo
$ws.Cells.Item($r, 7).Value = $productoId
$ws.Cells.Item($r, 8).Value = $producto
$ws.Cells.Item($r, 9).Value = $categoriaId
$ws.Cells.Item($r, 10).Value = $categoria
$ws.Cells.Item($r, 11).Value = "Start Ruby"
$ws.Cells.Item($r, 12).Value = "Primera"
$ws.Cells.Item($r, 13).Value = 80
$ws.Cells.Item($r, 14).Value = 14000
$ws.Cells.Item($r, 15).Value = 15000
$ws.Cells.Item($r, 16).Value = 14500
$ws.Cells.Item($r, 17).Value = $unidad
$ws.Cells.Item($r, 18).Value = $origen
$ws.Cells.Item($r, 19).Value = 1036
$ws.Cells.Item($r, 20).Value = $kgUnidad

# New row 331: Start Ruby / Segunda
$r = 331
$ws.Cells.Item($r, 1).Value = $mercadoId
$ws.Cells.Item($r, 2).Value = $mercado
$ws.Cells.Item($r, 3).Value = $region
$ws.Cells.Item($r, 4).Value = 44783
$ws.Cells.Item($r, 4).NumberFormat = $ws.Cells.Item(332, 4).NumberFormat
$ws.Cells.Item($r, 5).Value = $codreg
$ws.Cells.Item($r, 6).Value = $tipo
$ws.Cells.Item($r, 7).Value = $productoId
$ws.Cells.Item($r, 8).Value = $producto
$ws.Cells.Item($r, 9).Value = $categoriaId
$ws.Cells.Item($r, 10).Value = $categoria
$ws.Cells.Item($r, 11).Value = "Start Ruby"
$ws.Cells.Item($r, 12).Value = "Segunda"
$ws.Cells.Item($r, 13).Value = 40
$ws.Cells.Item($r, 14).Value = 12000
$ws.Cells.Item($r, 15).Value = 12000
$ws.Cells.Item($r, 16).Value = 12000
$ws.Cells.Item($r, 17).Value = $unidad
$ws.Cells.Item($r, 18).Value = $origen
$ws.Cells.Item($r, 19).Value = 857
$ws.Cells.Item($r, 20).Value = $kgUnidad
